# edit.ps1 — applies the commit's three changes to 自动阅读使用说明.docx
#  1) Insert a new paragraph ("！只需打开汇总帖即可，...") right after the
#     "登录微博网页，打开阅读汇总贴（建议将快捷方式放在桌面，比较好找）。" paragraph.
#  2) Reword "...就开始逐个点进微博正文页自动阅读了。..." to
#     "...就自动逐个点进微博正文页阅读了。" and split the run there so a
#     "_GoBack" bookmark sits right after "阅读了。" (before "自动阅读完成前...").
#     Because bookmark names are unique, re-adding "_GoBack" here also moves
#     it away from its old location later in the document.
#  3) (handled automatically by step 2 above) the old "_GoBack" bookmark
#     that used to sit after "时代表脚本在当前页面可以运行，否则说明脚本不可用。"
#     is removed, since Bookmarks.Add reuses/relocates the existing bookmark.

$d = $word.ActiveDocument

# --- Change 1: insert new paragraph after "登录微博网页..." -----------------
$rng = $d.Content
$rng.Find.Execute(
    "登录微博网页，打开阅读汇总贴（建议将快捷方式放在桌面，比较好找）。",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter("！只需打开汇总帖即可，不用自己点进去，也不要点进虎子发的博！")

# --- Change 2: reword the "打开开发者工具..." sentence ----------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "就开始逐个点进微博正文页自动阅读了。自动阅读完成前",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "就自动逐个点进微博正文页阅读了。自动阅读完成前", 2
) | Out-Null

# --- Change 3: move the "_GoBack" bookmark to right before "自动阅读完成前请" --
# (this both splits the run in two, as the diff shows, and removes the
#  bookmark from its old spot near "...脚本不可用。", since a bookmark name
#  can only exist once in the document.)
$rng3 = $d.Content
$rng3.Find.Execute(
    "自动阅读完成前请一直挂着页面",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$pt = $d.Range($rng3.Start, $rng3.Start)
$d.Bookmarks.Add("_GoBack", $pt) | Out-Null

Write-Output "ok"
